$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts NAME..COURSE from A..K to B..L)
$ws.Columns.Item(1).Insert()

# New "ID" header column (values intentionally left blank per the target data)
$ws.Range("A1").Value = "ID"

# New "ADDRESS" header + values in the new trailing column (M)
$ws.Range("M1").Value = "ADDRESS"

$address = "[ { street : 735 Olive Street }  ,  { city : Sunway City }  ,  { country : Zimbabwe } ]"
$ws.Range("M2").Value = $address
$ws.Range("M3").Value = $address

# Reorder the MARKS array values (2nd/3rd entries swapped) - now in column K
$ws.Range("K2").Value = "[[86.75, 90.45, 37.0]]"
$ws.Range("K3").Value = "[[57.0, 70.0, 56.05]]"

# Reorder the COURSE object list values (entries swapped) - now in column L
$ws.Range("L2").Value = "[{Maths : 2010} ,{Comp Science : 2012} ]"
$ws.Range("L3").Value = "[{Abuse : 2024} ,{Divorce : 2023} ]"
